$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three Registration Date values so they all land on 2026-03-01 (serial 46082)
$ws.Range("E2").Value = 46082
$ws.Range("E3").Value = 46082
$ws.Range("E4").Value = 46082

# Re-apply (refresh) the date number format across the existing date column/cells so the
# underlying number format record is rewritten from "m/d/yyyy"-style (numFmtId 14) to an
# explicit ISO format.
$ws.Range("E1:E4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Columns("E").NumberFormat = "yyyy\-mm\-dd;@"

# Add the new participant e-mail in D4, as a hyperlink
$ws.Range("D4").Value = "Bob@example.com"
[void]$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:Bob@example.com")
$ws.Range("D4").Style = "Normal"

# Move the active selection
$ws.Range("D7").Select() | Out-Null
